# Login functionality regression tests
# Adds a "forgot password" test-user row (row 6) to the "test" sheet,
# containing the demo account used for the forgot-password flow together
# with the most recently generated temporary password, styles it like the
# other rows (console/mono font for the generated values, Hyperlink style
# for the e-mail address), and leaves the "test" sheet as the active tab
# with the cursor parked below the new data (A7), matching the state
# after the regression-test run recorded the change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "test" sheet

# ---- New data row (row 6): forgot-password demo user ----------------
$ws.Cells.Item(6, 1).Value = "forgot password user"
$ws.Cells.Item(6, 2).Value = "qataskdemoaccnt@gmail.com"
$ws.Cells.Item(6, 3).Value = "newPT_910*811"

# Vertical-center the e-mail cell first so the hyperlink style (added next)
# keeps referencing the "Hyperlink" cell style (xfId 1) instead of "Normal".
$ws.Cells.Item(6, 2).VerticalAlignment = -4108   # xlVAlignCenter

# Hyperlink the e-mail address to its mailto: target (mirrors the other
# user/password rows above it).
$ws.Hyperlinks.Add($ws.Cells.Item(6, 2), "mailto:qataskdemoaccnt@gmail.com") | Out-Null

# Style the username / generated-password cells with the monospace
# "console output" look used for auto-generated test values.
foreach ($colIdx in @(1, 3)) {
    $cell = $ws.Cells.Item(6, $colIdx)
    $cell.Font.Name = "JetBrains Mono"
    $cell.Font.Size = 9.8
    $cell.Font.Color = 5867370   # BGR for RGB(0x6A,0x87,0x59) -> FF6A8759
    $cell.VerticalAlignment = -4108   # xlVAlignCenter
}

# Widen column A so the new "forgot password user" label fits (best-fit).
$ws.Columns.Item(1).ColumnWidth = 35.2

# ---- Active sheet / selection ----------------------------------------
# The "test" sheet becomes the active tab again (it was sheet "demo"
# before), with the selection resting on the row right below the new data.
$ws.Activate()
$ws.Range("A7").Select() | Out-Null
